$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractText = @'
Introduction
id="Par1">The coronavirus disease 2019 (COVID-19) pandemic has led to rapid adoption of teleophthalmology to deliver eyecare remotely.

 The purpose of our study was to assess the implementation and patient acceptability of video consultation for outpatient ophthalmic care at our institution.


Methods
id="Par2">We conducted a retrospective, cross-sectional analysis and patient survey of adult patients who completed a virtual video visit at our institution from 18 March 18 through to 27 April 2020. All video visit encounters were assessed for patient characteristics, diagnoses, management, and follow-up outcomes.

 Patients were surveyed for their feedback on acceptability and utility of their virtual video consultation.


Results
id="Par3">A total of 219 patients (mean age 55 years; range 21–89 years) completed 231 video visit encounters at our department over a 6-week period, of whom 118 were women (54%).

 About half of these encounters were acute visits (102 visits, 47%).

 The most common diagnosis of these visits was postoperative state (20 visits, 9% of the total), followed by conjunctivitis (16 visits, 7%), and keratitis (14 visits, 6%).

 The most common management decisions were medication prescription (102 visits, 46%) or reassurance (86 visits, 39%), while 17 video visit patients (8%) were escalated to an urgent, in-person evaluation.

 Ninety-two patients completed a follow-up survey (42% response rate), of whom 45 (49%) indicated that they might have delayed seeking care during this pandemic in the absence of a virtual video option.

 Seventy-two (78%) reported that they would consider participating in a video visit as an alternative to an office-based encounter in the future, and the overall video visit experience was rated highly, with a weighted mean Likert scale rating of 4.3 out of 5 (Cronbach’s α = 0.88).


Conclusion
id="Par4">Virtual video visits may be used to manage a range of ophthalmic complaints.

 Patients participating in this survey found such video visits acceptable and timesaving, and the majority would consider using video consultations for future eyecare encounters.



'@

$authorsText = @'
[Gagan%Kalra%NULL%1,   Andrew M.%Williams%NULL%1,   Patrick W.%Commiskey%NULL%1,   Eve M. R.%Bowers%NULL%1,   Tadhg%Schempf%NULL%1,   José-Alain%Sahel%NULL%1,   Evan L.%Waxman%waxmane@upmc.edu%1,   Roxana%Fu%fur3@upmc.edu%1]
'@

$ws.Range("D3").Value = $abstractText
$ws.Range("E3").Value = $authorsText
